# EarGear Protocol v1.docx edit script
# Applies:
#  1. "Returns firmware version" -> "Returns " / "VER plus " / "firmware version" (3 runs)
#  2. "Enter Listen Mode" paragraph gains a trailing ", returns LISTEN ON" run
#  3. "Stop Listen Mode" paragraph gains a trailing ", returns LISTEN OFF" run
#  4. "Returns "OK"" -> "Returns "PONG"" with the _GoBack bookmark now sitting
#     between the PONG run and the closing-quote run
#  5. The old _GoBack bookmark (previously after "DSSP") is removed

$d = $word.ActiveDocument

# Helper: insert `$text` at the collapsed position `$insertPos`, then apply the
# document's body font to the freshly-inserted run so it matches the existing
# "Century Gothic" formatting used throughout this document.
function Insert-FormattedRun($insertPos, $text) {
    $ins = $d.Range($insertPos, $insertPos)
    $ins.InsertAfter($text)
    $newRange = $d.Range($insertPos, $insertPos + $text.Length)
    $newRange.Font.Name = "Century Gothic"
    return $insertPos + $text.Length
}

# Helper: find the first paragraph whose text starts with `$prefix`.
function Find-Paragraph($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($prefix)) {
            return $d.Paragraphs($i)
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Edit 1: VER command description - split "Returns firmware version" into
# "Returns " + "VER plus " + "firmware version". The run holding this text
# also carries a leading <w:tab/> sibling element, so deleting/replacing text
# that begins right at the run's start would flatten that <w:tab/> into a
# literal tab character when the run gets re-serialized. Instead, we use a
# throwaway Bookmarks.Add/Delete pair to cleanly carve "firmware version"
# (the part being removed) into its own run first -- Bookmarks.Add splits
# runs without collapsing sibling <w:tab/> elements -- leaving "Returns "
# (plus the original tab) untouched in the original run, exactly as the
# target markup expects.
# ---------------------------------------------------------------------------
$p = Find-Paragraph "VER"
$full = $p.Range.Text
$idx = $full.IndexOf("firmware version")
$start = $p.Range.Start + $idx
$end = $p.Range.End - 1

$splitPoint = $d.Range($start, $start)
$d.Bookmarks.Add("_TempSplit1", $splitPoint) | Out-Null
$d.Range($start, $end).Text = ""
$d.Bookmarks.Item("_TempSplit1").Delete()

$pos = $start
$pos = Insert-FormattedRun $pos "VER plus "
$pos = Insert-FormattedRun $pos "firmware version"

# ---------------------------------------------------------------------------
# Edit 2: LISTEN command description gains ", returns LISTEN ON"
# ---------------------------------------------------------------------------
$p = Find-Paragraph "LISTEN"
$endOfText = $p.Range.End - 1
Insert-FormattedRun $endOfText ", returns LISTEN ON" | Out-Null

# ---------------------------------------------------------------------------
# Edit 3: ENDLISTEN command description gains ", returns LISTEN OFF"
# ---------------------------------------------------------------------------
$p = Find-Paragraph "ENDLISTEN"
$endOfText = $p.Range.End - 1
Insert-FormattedRun $endOfText ", returns LISTEN OFF" | Out-Null

# ---------------------------------------------------------------------------
# Edit 5 (done first): remove the old _GoBack bookmark that currently sits
# right after "DSSP", since a document can only contain one bookmark named
# "_GoBack" at a time (it is Word's hidden "last edit" marker) and the diff
# relocates it down to the PING paragraph below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Edit 4: PING command description - "Returns "OK"" -> "Returns "PONG"", with
# the _GoBack bookmark now sitting right after "PONG" (before the closing
# curly quote). The run holding this text also carries a leading <w:tab/>
# element; splitting the run via a plain text/range edit would collapse that
# <w:tab/> into a literal tab character, so we first use Bookmarks.Add
# (which splits runs without flattening sibling <w:tab/> elements) to carve
# "OK" into its own run, then rewrite that isolated run's text normally.
# ---------------------------------------------------------------------------
$p = Find-Paragraph "PING"
$full = $p.Range.Text
$idx = $full.IndexOf("OK")
$start = $p.Range.Start + $idx
$end = $start + "OK".Length

$splitPoint = $d.Range($start, $start)
$d.Bookmarks.Add("_TempSplit", $splitPoint) | Out-Null

$d.Range($start, $end).Text = "PONG"
$d.Bookmarks.Item("_TempSplit").Delete()

$pongEnd = $start + "PONG".Length
$gobackRange = $d.Range($pongEnd, $pongEnd)
$d.Bookmarks.Add("_GoBack", $gobackRange) | Out-Null

Write-Output "VER:       [$((Find-Paragraph 'VER').Range.Text)]"
Write-Output "LISTEN:    [$((Find-Paragraph 'LISTEN').Range.Text)]"
Write-Output "ENDLISTEN: [$((Find-Paragraph 'ENDLISTEN').Range.Text)]"
Write-Output "PING:      [$((Find-Paragraph 'PING').Range.Text)]"
Write-Output "GoBack exists: $($d.Bookmarks.Exists('_GoBack'))"
